# Apply odds updates to Sheet1 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 (Beijing Guoan vs Qingdao West Coast) ----
$ws.Cells.Item(3, 7).Value  = 1.33   # G3
$ws.Cells.Item(3, 9).Value  = 7.5    # I3
$ws.Cells.Item(3, 10).Value = 1.02   # J3
$ws.Cells.Item(3, 11).Value = 21     # K3
$ws.Cells.Item(3, 12).Value = 1.13   # L3
$ws.Cells.Item(3, 13).Value = 6      # M3
$ws.Cells.Item(3, 14).Value = 1.44   # N3
$ws.Cells.Item(3, 15).Value = 2.7    # O3
$ws.Cells.Item(3, 16).Value = 1.22   # P3
$ws.Cells.Item(3, 17).Value = 4      # Q3
$ws.Cells.Item(3, 18).Value = 1.75   # R3
$ws.Cells.Item(3, 19).Value = 2      # S3
$ws.Cells.Item(3, 20).Value = 9.5    # T3
$ws.Cells.Item(3, 21).Value = 8      # U3
$ws.Cells.Item(3, 23).Value = 9.5    # W3
$ws.Cells.Item(3, 25).Value = 21     # Y3
$ws.Cells.Item(3, 26).Value = 21     # Z3
$ws.Cells.Item(3, 28).Value = 19     # AB3
$ws.Cells.Item(3, 29).Value = 51     # AC3
$ws.Cells.Item(3, 30).Value = 201    # AD3
$ws.Cells.Item(3, 31).Value = 23     # AE3
$ws.Cells.Item(3, 35).Value = 51     # AI3

# ---- Row 5 (Metta vs BFC Daugavpils) : previously blank cells now filled ----
$ws.Cells.Item(5, 7).Value  = 3      # G5
$ws.Cells.Item(5, 8).Value  = 3.45   # H5
$ws.Cells.Item(5, 9).Value  = 2.07   # I5
$ws.Cells.Item(5, 14).Value = 1.8    # N5
$ws.Cells.Item(5, 15).Value = 1.8    # O5
$ws.Cells.Item(5, 20).Value = 8.25   # T5
$ws.Cells.Item(5, 21).Value = 13     # U5
$ws.Cells.Item(5, 22).Value = 9.25   # V5
$ws.Cells.Item(5, 23).Value = 29     # W5
$ws.Cells.Item(5, 24).Value = 20     # X5
$ws.Cells.Item(5, 25).Value = 26     # Y5
$ws.Cells.Item(5, 26).Value = 10.75  # Z5
$ws.Cells.Item(5, 27).Value = 5.9    # AA5
$ws.Cells.Item(5, 28).Value = 11.75  # AB5
$ws.Cells.Item(5, 29).Value = 45     # AC5
$ws.Cells.Item(5, 30).Value = 300    # AD5
$ws.Cells.Item(5, 31).Value = 6.7    # AE5
$ws.Cells.Item(5, 32).Value = 8.5    # AF5
$ws.Cells.Item(5, 33).Value = 7.5    # AG5
$ws.Cells.Item(5, 34).Value = 15.5   # AH5
$ws.Cells.Item(5, 35).Value = 13.5   # AI5
$ws.Cells.Item(5, 36).Value = 21     # AJ5

# ---- Row 6 (Tukums 2000 vs Jelgava) : previously blank cells now filled ----
$ws.Cells.Item(6, 7).Value  = 3.2    # G6
$ws.Cells.Item(6, 8).Value  = 3.25   # H6
$ws.Cells.Item(6, 9).Value  = 2.05   # I6
$ws.Cells.Item(6, 14).Value = 1.82   # N6
$ws.Cells.Item(6, 15).Value = 1.78   # O6
$ws.Cells.Item(6, 20).Value = 8.5    # T6
$ws.Cells.Item(6, 21).Value = 14     # U6
$ws.Cells.Item(6, 22).Value = 9.5    # V6
$ws.Cells.Item(6, 23).Value = 32     # W6
$ws.Cells.Item(6, 24).Value = 22     # X6
$ws.Cells.Item(6, 25).Value = 27     # Y6
$ws.Cells.Item(6, 26).Value = 10     # Z6
$ws.Cells.Item(6, 27).Value = 5.6    # AA6
$ws.Cells.Item(6, 28).Value = 11.25  # AB6
$ws.Cells.Item(6, 29).Value = 45     # AC6
$ws.Cells.Item(6, 30).Value = 300    # AD6
$ws.Cells.Item(6, 31).Value = 6.7    # AE6
$ws.Cells.Item(6, 32).Value = 8.5    # AF6
$ws.Cells.Item(6, 33).Value = 7.3    # AG6
$ws.Cells.Item(6, 34).Value = 15.5   # AH6
$ws.Cells.Item(6, 35).Value = 13     # AI6
$ws.Cells.Item(6, 36).Value = 20     # AJ6

# ---- Row 7 (Super Nova vs FK Liepaja) : previously blank cells now filled ----
$ws.Cells.Item(7, 7).Value  = 2.47   # G7
$ws.Cells.Item(7, 8).Value  = 3.35   # H7
$ws.Cells.Item(7, 9).Value  = 2.45   # I7
$ws.Cells.Item(7, 14).Value = 1.82   # N7
$ws.Cells.Item(7, 15).Value = 1.78   # O7
$ws.Cells.Item(7, 16).Value = 1.37   # P7
$ws.Cells.Item(7, 17).Value = 2.5    # Q7
$ws.Cells.Item(7, 20).Value = 7.8    # T7
$ws.Cells.Item(7, 21).Value = 11.75  # U7
$ws.Cells.Item(7, 22).Value = 8.75   # V7
$ws.Cells.Item(7, 23).Value = 25     # W7
$ws.Cells.Item(7, 24).Value = 18.5   # X7
$ws.Cells.Item(7, 25).Value = 25     # Y7
$ws.Cells.Item(7, 26).Value = 10     # Z7
$ws.Cells.Item(7, 27).Value = 5.7    # AA7
$ws.Cells.Item(7, 28).Value = 11.5   # AB7
$ws.Cells.Item(7, 29).Value = 45     # AC7
$ws.Cells.Item(7, 30).Value = 300    # AD7
$ws.Cells.Item(7, 31).Value = 7      # AE7
$ws.Cells.Item(7, 32).Value = 9.5    # AF7
$ws.Cells.Item(7, 33).Value = 7.7    # AG7
$ws.Cells.Item(7, 34).Value = 18     # AH7
$ws.Cells.Item(7, 35).Value = 14.5   # AI7
$ws.Cells.Item(7, 36).Value = 22     # AJ7

# ---- Row 8 (Daegu vs Pohang) ----
$ws.Cells.Item(8, 14).Value = 1.88   # N8
$ws.Cells.Item(8, 15).Value = 1.93   # O8

# ---- Row 9 (Jeonbuk vs Suwon FC) ----
$ws.Cells.Item(9, 7).Value  = 1.6    # G9
$ws.Cells.Item(9, 8).Value  = 3.9    # H9
$ws.Cells.Item(9, 9).Value  = 5.5    # I9
$ws.Cells.Item(9, 20).Value = 7.5    # T9
$ws.Cells.Item(9, 23).Value = 12     # W9

# ---- Row 10 (Seoul vs Gangwon) ----
$ws.Cells.Item(10, 7).Value  = 1.67  # G10
$ws.Cells.Item(10, 8).Value  = 3.6   # H10
$ws.Cells.Item(10, 10).Value = 1.05  # J10
$ws.Cells.Item(10, 12).Value = 1.33  # L10
$ws.Cells.Item(10, 20).Value = 5.5   # T10
$ws.Cells.Item(10, 23).Value = 12    # W10
$ws.Cells.Item(10, 28).Value = 21    # AB10
$ws.Cells.Item(10, 31).Value = 11    # AE10

$wb.Save()
